$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (46075 -> 46076) for every data row (rows 2 through 371).
for ($row = 2; $row -le 371; $row++) {
    $ws.Cells.Item($row, 3).Value = 46076
}
